$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 223, shifting existing rows 223-319 down to 224-320
$ws.Rows.Item(223).Insert()

# Populate the newly inserted row 223 with the new record
$ws.Cells.Item(223, 1).Value = 4
$ws.Cells.Item(223, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(223, 3).Value = "Los Lagos"
$ws.Cells.Item(223, 4).Value = 44755
$ws.Cells.Item(223, 5).Value = 10
$ws.Cells.Item(223, 6).Value = 100114014
$ws.Cells.Item(223, 7).Value = "Betarraga"
$ws.Cells.Item(223, 8).Value = "Sin especificar"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 250
$ws.Cells.Item(223, 11).Value = 1200
$ws.Cells.Item(223, 12).Value = 1200
$ws.Cells.Item(223, 13).Value = 1200
$ws.Cells.Item(223, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(223, 15).Value = "Región del Maule"
$ws.Cells.Item(223, 16).Value = 240
$ws.Cells.Item(223, 17).Value = 5
$ws.Cells.Item(223, 18).Value = "Hortaliza"
